# NIT-9008924770.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# Data update on sheet "Hoja1":
#  - Rows 16 and 17 swap their "Periodo Mora" (col E) / "Valor Mora" (col F)
#    values (period 1611 <-> 1608, with their matching mora amounts).
#  - Row 18's "Salario Basico" (col G) is corrected from 1,000,000 to 908,526.
#  - The company logo picture is nudged 13.5pt (171450 EMU) to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Swap the two "Periodo Mora" / "Valor Mora" data rows -------------------
$ws.Range("E16").Value = "1608"
$ws.Range("F16").Value = 2758

$ws.Range("E17").Value = "1611"
$ws.Range("F17").Value = 3677

# --- Correct the "Salario Basico" figure on the last data row --------------
$ws.Range("G18").Value = 908526

# --- Reposition the logo image: shift left by 171450 EMU (13.5 pt) ---------
# Recompute the absolute left position precisely (in EMU) from the current
# width of column A (in points) plus the picture's stored horizontal offset
# into column B, so the move lands on an exact EMU value instead of drifting
# from the 2-decimal rounding that the .Left getter applies.
$shp = $ws.Shapes.Item(1)

$colAWidthEMU = $ws.Columns.Item(1).Width * 12700.0
$currentColOffEMU = 667900.0
$targetColOffEMU = 496450.0

$currentLeftEMU = $colAWidthEMU + $currentColOffEMU
$newLeftEMU = $currentLeftEMU + ($targetColOffEMU - $currentColOffEMU)

$shp.Left = $newLeftEMU / 12700.0
